$wb = $excel.ActiveWorkbook
$frSheet = $wb.Worksheets.Item("FR_Properties")
Write-Host "FR_Properties index: $($frSheet.Index)"
$newSheet = $wb.Worksheets.Add($frSheet)
$newSheet.Name = "Grouper resultat"
Write-Host "New sheet index: $($newSheet.Index)"
Write-Host "New sheet name: $($newSheet.Name)"
